$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update shared-string header text (volume number and report date range) ---
$a8 = $ws.Range("A8")
$a8Text = $a8.Value2
$a8Pos = $a8Text.LastIndexOf("5") + 1
$a8.Characters($a8Pos, 1).Text = "6"

$c9 = $ws.Range("C9")
$c9Text = $c9.Value2
$c9Pos1 = $c9Text.IndexOf("1/30/2023") + 1
$c9.Characters($c9Pos1, 9).Text = "2/6/2023"
$c9Text2 = $c9.Value2
$c9Pos2 = $c9Text2.IndexOf("2/5/2023") + 1
$c9.Characters($c9Pos2, 8).Text = "2/12/2023"

# --- Update crime-statistics data cells (rows 15-30) ---

# Row 15
$ws.Range("C15").NumberFormat = '#,##0'
$ws.Range("C15").Value = 1
$ws.Range("F15").Value = 2
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 100
$ws.Range("I15").Value = 4
$ws.Range("K15").Value = 100
$ws.Range("L15").Value = 33.333333333333
$ws.Range("N15").Value = -42.857142857142

# Row 16
$ws.Range("C16").Value = 8
$ws.Range("E16").Value = 700
$ws.Range("F16").Value = 18
$ws.Range("G16").Value = 8
$ws.Range("H16").Value = 125
$ws.Range("I16").Value = 28
$ws.Range("J16").Value = 11
$ws.Range("K16").Value = 154.545454545455
$ws.Range("L16").Value = 154.545454545455
$ws.Range("M16").Value = -12.5
$ws.Range("N16").Value = -80.952380952380

# Row 17
$ws.Range("C17").Value = 8
$ws.Range("E17").Value = 166.666666666667
$ws.Range("F17").Value = 20
$ws.Range("G17").Value = 16
$ws.Range("H17").Value = 25
$ws.Range("I17").Value = 31
$ws.Range("J17").Value = 23
$ws.Range("K17").Value = 34.782608695652
$ws.Range("L17").Value = 106.666666666667
$ws.Range("M17").Value = 93.75
$ws.Range("N17").Value = -43.636363636363

# Row 18
$ws.Range("C18").Value = 4
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 16
$ws.Range("G18").Value = 19
$ws.Range("H18").Value = -15.789473684210
$ws.Range("I18").Value = 32
$ws.Range("J18").Value = 42
$ws.Range("K18").Value = -23.809523809523
$ws.Range("L18").Value = 18.518518518518
$ws.Range("M18").Value = -41.818181818181
$ws.Range("N18").Value = -76.470588235294

# Row 19
$ws.Range("C19").Value = 19
$ws.Range("D19").Value = 15
$ws.Range("E19").Value = 26.666666666666
$ws.Range("F19").Value = 76
$ws.Range("G19").Value = 46
$ws.Range("H19").Value = 65.217391304347
$ws.Range("I19").Value = 104
$ws.Range("J19").Value = 64
$ws.Range("K19").Value = 62.5
$ws.Range("L19").Value = 147.619047619048
$ws.Range("M19").Value = 136.363636363636
$ws.Range("N19").Value = 96.226415094339

# Row 20
$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 5
$ws.Range("E20").Value = -20
$ws.Range("I20").Value = 16
$ws.Range("J20").Value = 19
$ws.Range("K20").Value = -15.789473684210
$ws.Range("L20").Value = 14.285714285714
$ws.Range("M20").Value = -5.882352941176
$ws.Range("N20").Value = -83.333333333333

# Row 21
$ws.Range("C21").Value = 44
$ws.Range("D21").Value = 28
$ws.Range("E21").Value = 57.142857142857
$ws.Range("F21").Value = 141
$ws.Range("G21").Value = 100
$ws.Range("H21").Value = 41
$ws.Range("I21").Value = 216
$ws.Range("J21").Value = 161
$ws.Range("K21").Value = 34.161490683229
$ws.Range("L21").Value = 92.857142857142
$ws.Range("M21").Value = 30.909090909090
$ws.Range("N21").Value = -56.451612903225

# Row 22
$ws.Range("C22").NumberFormat = '#,##0'
$ws.Range("C22").Value = 3
$ws.Range("D22").NumberFormat = '#,##0'
$ws.Range("D22").Value = 1
$ws.Range("E22").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E22").Value = 200
$ws.Range("F22").Value = 3
$ws.Range("H22").Value = -25
$ws.Range("I22").Value = 5
$ws.Range("J22").Value = 5
$ws.Range("K22").Value = 0
$ws.Range("M22").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("M22").Value = 150

# Row 23
$ws.Range("C23").Value = 4
$ws.Range("D23").Value = 7
$ws.Range("E23").Value = -42.857142857142
$ws.Range("F23").Value = 8
$ws.Range("G23").Value = 14
$ws.Range("H23").Value = -42.857142857142
$ws.Range("I23").Value = 18
$ws.Range("J23").Value = 20
$ws.Range("K23").Value = -10
$ws.Range("L23").Value = 38.461538461538
$ws.Range("M23").Value = 38.461538461538

# Row 24
$ws.Range("C24").Value = 17
$ws.Range("D24").Value = 13
$ws.Range("E24").Value = 30.769230769230
$ws.Range("F24").Value = 90
$ws.Range("G24").Value = 83
$ws.Range("H24").Value = 8.433734939759
$ws.Range("I24").Value = 132
$ws.Range("J24").Value = 120
$ws.Range("K24").Value = 10
$ws.Range("L24").Value = 36.082474226804
$ws.Range("M24").Value = 0.763358778625

# Row 25
$ws.Range("C25").Value = 12
$ws.Range("D25").Value = 3
$ws.Range("E25").Value = 300
$ws.Range("F25").Value = 39
$ws.Range("G25").Value = 31
$ws.Range("H25").Value = 25.806451612903
$ws.Range("I25").Value = 55
$ws.Range("J25").Value = 42
$ws.Range("K25").Value = 30.952380952381
$ws.Range("L25").Value = 57.142857142857
$ws.Range("M25").Value = 14.583333333333

# Row 26
$ws.Range("C26").Value = 1
$ws.Range("F26").Value = 4
$ws.Range("G26").Value = 1
$ws.Range("H26").Value = 300
$ws.Range("I26").Value = 7
$ws.Range("K26").Value = 250
$ws.Range("L26").Value = 75

# Row 27
$ws.Range("C27").NumberFormat = '#,##0'
$ws.Range("C27").Value = 1
$ws.Range("D27").NumberFormat = '#,##0'
$ws.Range("D27").Value = 2
$ws.Range("E27").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E27").Value = -50
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 4
$ws.Range("J27").Value = 5
$ws.Range("K27").Value = -20
$ws.Range("L27").Value = -50

# Row 28
$ws.Range("C28").NumberFormat = '#,##0'
$ws.Range("C28").Value = 5
$ws.Range("D28").NumberFormat = '#,##0'
$ws.Range("D28").Value = 1
$ws.Range("E28").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E28").Value = 400
$ws.Range("F28").NumberFormat = '#,##0'
$ws.Range("F28").Value = 5
$ws.Range("G28").Value = 2
$ws.Range("H28").Value = 150
$ws.Range("I28").NumberFormat = '#,##0'
$ws.Range("I28").Value = 5
$ws.Range("J28").Value = 2
$ws.Range("K28").Value = 150
$ws.Range("M28").Value = 400
$ws.Range("N28").Value = -64.285714285714

# Row 29
$ws.Range("C29").NumberFormat = '#,##0'
$ws.Range("C29").Value = 2
$ws.Range("D29").NumberFormat = '#,##0'
$ws.Range("D29").Value = 1
$ws.Range("E29").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E29").Value = 100
$ws.Range("F29").NumberFormat = '#,##0'
$ws.Range("F29").Value = 2
$ws.Range("G29").Value = 2
$ws.Range("H29").Value = 0
$ws.Range("I29").NumberFormat = '#,##0'
$ws.Range("I29").Value = 2
$ws.Range("J29").Value = 2
$ws.Range("K29").Value = 0
$ws.Range("M29").Value = 100
$ws.Range("N29").Value = -85.714285714285

# Row 30
$ws.Range("D30").Value = 3
$ws.Range("G30").Value = 4
$ws.Range("J30").Value = 5
